$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -5
    3  = -2
    4  = -4
    5  = 0
    6  = 1
    7  = -3
    8  = 1
    9  = 7
    10 = -2
    11 = 1
    12 = 8
    13 = -1
    14 = -5
    15 = 0
    16 = -3
    17 = 1
    20 = 1
    22 = -1
    23 = -2
    24 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
